$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete 2004/2008/2009 data rows (old rows 2:4); this shifts
# 2010..2020 up to rows 2..12, matching the new layout.
$ws.Range("A2:A4").EntireRow.Delete()

# Append the new 2021 data row as row 13, copying the year-label cell's
# format (bold font + border + centered alignment) from the row above it.
$ws.Cells.Item(12,1).Copy()
$ws.Cells.Item(13,1).PasteSpecial(-4122)
$ws.Cells.Item(13,1).Value = "2021年"

$ws.Cells.Item(13,2).Value = 144496
$ws.Cells.Item(13,3).Value = 186192
$ws.Cells.Item(13,4).Value = 43580
$ws.Cells.Item(13,5).Value = 913223
$ws.Cells.Item(13,6).Value = 264751
$ws.Cells.Item(13,7).Value = 248716
$ws.Cells.Item(13,8).Value = 132702
$ws.Cells.Item(13,9).Value = 283805
$ws.Cells.Item(13,10).Value = 1155744
$ws.Cells.Item(13,11).Value = 205738
$ws.Cells.Item(13,12).Value = 249212
$ws.Cells.Item(13,13).Value = 159136
$ws.Cells.Item(13,14).Value = 291764
$ws.Cells.Item(13,15).Value = 1499459
$ws.Cells.Item(13,16).Value = 126518
$ws.Cells.Item(13,17).Value = 976345
$ws.Cells.Item(13,18).Value = 118438
$ws.Cells.Item(13,19).Value = 31124
$ws.Cells.Item(13,20).Value = 375895
$ws.Cells.Item(13,21).Value = 1198973
